$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.974.26"
$ws.Range("E2").Value = "  -2.35%  "
$ws.Range("D3").Value = "2.498.69"
$ws.Range("E3").Value = "  -3.77%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "550.70"
$ws.Range("E5").Value = "  -3.80%  "
$ws.Range("D6").Value = "146.89"
$ws.Range("E6").Value = "  -5.05%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "0.614"
$ws.Range("E8").Value = "  -1.29%  "
$ws.Range("D9").Value = "2.498.07"
$ws.Range("E9").Value = "  -3.74%  "
$ws.Range("E10").Value = "  -9.40%  "
$ws.Range("E11").Value = "  -1.61%  "
$ws.Range("E12").Value = "  -7.99%  "
$ws.Range("E13").Value = "  -6.22%  "
$ws.Range("D14").Value = "26.13"
$ws.Range("E14").Value = "  -6.84%  "
$ws.Range("D15").Value = "2.948.98"
$ws.Range("E15").Value = "  -3.63%  "
$ws.Range("D16").Value = "61.833.48"
$ws.Range("E16").Value = "  -2.29%  "
$ws.Range("E17").Value = "  -8.25%  "
$ws.Range("D18").Value = "2.495.06"
$ws.Range("E18").Value = "  -3.47%  "
$ws.Range("E19").Value = "  -7.24%  "
$ws.Range("D20").Value = "7.02"
$ws.Range("E20").Value = "  -6.09%  "
$ws.Range("E21").Value = "  -7.57%  "
$ws.Range("D22").Value = "321.49"
$ws.Range("E22").Value = "  -5.79%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "63.86"
$ws.Range("E24").Value = "  -5.18%  "
$ws.Range("E25").Value = "  -4.06%  "
$ws.Range("E26").Value = "  -6.33%  "
$ws.Range("D27").Value = "2.623.22"
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("E29").Value = "  -5.12%  "
$ws.Range("D30").Value = "8.39"
$ws.Range("E30").Value = "  -8.13%  "
$ws.Range("D31").Value = "534.29"
$ws.Range("E31").Value = "  -7.41%  "
$ws.Range("E32").Value = "  -3.25%  "
$ws.Range("E33").Value = "  -6.69%  "
$ws.Range("E34").Value = "  -7.80%  "
$ws.Range("E35").Value = "  -8.96%  "
$ws.Range("E36").Value = "  -10.15%  "
$ws.Range("D37").Value = "4.86"
$ws.Range("E37").Value = "  -10.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("E39").Value = "  -5.85%  "
$ws.Range("D40").Value = "18.49"
$ws.Range("E40").Value = "  -6.17%  "
$ws.Range("D41").Value = "143.99"
$ws.Range("E41").Value = "  -6.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("E43").Value = "  -8.99%  "
$ws.Range("D44").Value = "40.36"
$ws.Range("E44").Value = "  -2.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.30"
$ws.Range("E45").Value = "  -7.42%  "
$ws.Range("D46").Value = "148.73"
$ws.Range("E46").Value = "  -5.01%  "
$ws.Range("E47").Value = "  -8.19%  "
$ws.Range("D48").Value = "20.76"
$ws.Range("E48").Value = "  -10.55%  "
$ws.Range("E49").Value = "  -8.69%  "
$ws.Range("E50").Value = "  -5.93%  "
$ws.Range("D51").Value = "0.0948"
$ws.Range("E51").Value = "  -5.20%  "
